$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Trening" column header (F1), reusing the existing header style ---
$ws.Range("F1").Value = "Trening"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "Trening"

# --- Custom date/time display format for column A (rows 2-13) ---
# First assignment registers numFmt 164 ("yyyy-mm-dd h:mm:ss") on A2 and
# is immediately replaced by numFmt 165 ("YYYY-MM-DD HH:MM:SS"), matching the
# two <numFmt> entries kept in styles.xml (164 stays registered but unused).
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
for ($r = 3; $r -le 13; $r++) {
    $ws.Range("A$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# --- Update existing rows 2-7 and append new rows 8-13 ---
$ws.Range("A2").Value = 45684.59094606482
$ws.Range("B2").Value = 457.7
$ws.Range("C2").Value = 10.62
$ws.Range("D2").Value = 2.356715134211947
$ws.Range("E2").Value = "10-15"
$ws.Range("F2").Value = "Duża Gra"

$ws.Range("A3").Value = 45684.59206180555
$ws.Range("B3").Value = 554.1
$ws.Range("C3").Value = 10.3
$ws.Range("D3").Value = 2.037503021103996
$ws.Range("E3").Value = "10-15"
$ws.Range("F3").Value = "Duża Gra"

$ws.Range("A4").Value = 45684.59324351852
$ws.Range("B4").Value = 656.2
$ws.Range("C4").Value = 10.01
$ws.Range("D4").Value = 1.827428545270647
$ws.Range("E4").Value = "10-15"
$ws.Range("F4").Value = "Duża Gra"

$ws.Range("A5").Value = 45684.59094375
$ws.Range("B5").Value = 457.5
$ws.Range("C5").Value = 8.98
$ws.Range("D5").Value = 2.416621276310511
$ws.Range("E5").Value = "5-10"
$ws.Range("F5").Value = "Duża Gra"

$ws.Range("A6").Value = 45684.59206064814
$ws.Range("B6").Value = 554
$ws.Range("C6").Value = 9.51
$ws.Range("D6").Value = 2.041833485875812
$ws.Range("E6").Value = "5-10"
$ws.Range("F6").Value = "Duża Gra"

$ws.Range("A7").Value = 45684.59324236111
$ws.Range("B7").Value = 656.1
$ws.Range("C7").Value = 9.41
$ws.Range("D7").Value = 1.844880853380476
$ws.Range("E7").Value = "5-10"
$ws.Range("F7").Value = "Duża Gra"

$ws.Range("A8").Value = 45684.59792986111
$ws.Range("B8").Value = 1061.1
$ws.Range("C8").Value = 14.03
$ws.Range("D8").Value = 3.011687414986746
$ws.Range("E8").Value = "10-15"
$ws.Range("F8").Value = "Mała Gra"

$ws.Range("A9").Value = 45684.60068101852
$ws.Range("B9").Value = 1298.8
$ws.Range("C9").Value = 14.66
$ws.Range("D9").Value = 2.598737512316023
$ws.Range("E9").Value = "10-15"
$ws.Range("F9").Value = "Mała Gra"

$ws.Range("A10").Value = 45684.60428171296
$ws.Range("B10").Value = 1609.9
$ws.Range("C10").Value = 13.54
$ws.Range("D10").Value = 2.660114662987844
$ws.Range("E10").Value = "10-15"
$ws.Range("F10").Value = "Mała Gra"

$ws.Range("A11").Value = 45684.60065787037
$ws.Range("B11").Value = 1296.8
$ws.Range("C11").Value = 7.52
$ws.Range("D11").Value = 2.329664622034346
$ws.Range("E11").Value = "5-10"
$ws.Range("F11").Value = "Mała Gra"

$ws.Range("A12").Value = 45684.60141134259
$ws.Range("B12").Value = 1361.9
$ws.Range("C12").Value = 8.37
$ws.Range("D12").Value = 2.29239056791578
$ws.Range("E12").Value = "5-10"
$ws.Range("F12").Value = "Mała Gra"

$ws.Range("A13").Value = 45684.60220416667
$ws.Range("B13").Value = 1430.4
$ws.Range("C13").Value = 9.27
$ws.Range("D13").Value = 2.648619515555246
$ws.Range("E13").Value = "5-10"
$ws.Range("F13").Value = "Mała Gra"
